$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login Page")
$vehicleSheet = $wb.Worksheets.Item("Selected Vehicle Page")

# Add the new data row on the "Login Page" sheet (leading apostrophe keeps it text,
# matching the existing text-formatted numbers in that column)
$loginSheet.Range("A4").Value = "'8800996794"

# Update the active cell / selection on the "Selected Vehicle Page" sheet
$vehicleSheet.Select()
$vehicleSheet.Range("A2").Select()

# Make "Login Page" the active sheet and set its selection, so it becomes the
# tab that is active/selected when the workbook is next opened
$loginSheet.Select()
$loginSheet.Range("A5").Select()
